# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to reflect the latest generated data (commit: output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 7294
    $ws.Range("F4").Value = 5543
    $ws.Range("F5").Value = 81
    $ws.Range("F14").Value = 646
    $ws.Range("F15").Value = 310
    $ws.Range("F19").Value = 43
}
